$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.427.28"
$ws.Range("E2").Value = "  +3.13%  "

$ws.Range("D3").Value = "3.368.37"
$ws.Range("E3").Value = "  +4.60%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'191.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.51%  "

$ws.Range("D6").Value = "'592.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.52%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").Value = "'0.134"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.09%  "

$ws.Range("D10").Value = "'6.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.11%  "

$ws.Range("E11").Value = "  +2.71%  "

$ws.Range("D12").Value = "3.957.47"
$ws.Range("E12").Value = "  +4.83%  "

$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").Value = "'28.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.69%  "

$ws.Range("D15").Value = "69.501.16"
$ws.Range("E15").Value = "  +3.13%  "

$ws.Range("E16").Value = "  +2.28%  "

$ws.Range("D17").Value = "3.347.76"
$ws.Range("E17").Value = "  +4.62%  "

$ws.Range("D18").Value = "'450.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +14.11%  "

$ws.Range("E19").Value = "  +1.87%  "

$ws.Range("D20").Value = "'13.84"
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = "  +3.91%  "

$ws.Range("D22").Value = "'74.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.93%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "3.517.67"
$ws.Range("E24").Value = "  +4.64%  "

$ws.Range("E25").Value = "  +4.85%  "

$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("E27").Value = "  +3.95%  "

$ws.Range("D28").Value = "'9.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.21%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").Value = "'1.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.33%  "

$ws.Range("D31").Value = "'23.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.56%  "

$ws.Range("E32").Value = "  +2.15%  "

$ws.Range("D33").Value = "'1.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.13%  "

$ws.Range("D34").Value = "'7.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  +4.90%  "

$ws.Range("D37").Value = "'165.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.15%  "

$ws.Range("E38").Value = "  +3.29%  "

$ws.Range("D39").Value = "'27.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.19%  "

$ws.Range("D40").Value = "'0.817"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.12%  "

$ws.Range("D41").Value = "'4.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.95%  "

$ws.Range("E42").Value = "  +1.02%  "

$ws.Range("D43").Value = "2.737.83"
$ws.Range("E43").Value = "  +5.72%  "

$ws.Range("E44").Value = "  +3.89%  "

$ws.Range("D45").Value = "'25.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.11%  "

$ws.Range("D46").Value = "'0.0690"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.28%  "

$ws.Range("D47").Value = "'343.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.20%  "

$ws.Range("D48").Value = "'40.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").Value = "'0.0285"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.14%  "

$ws.Range("D50").Value = "'32.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.26%  "

$ws.Range("E51").Value = "  +5.71%  "
